$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump "Forandrad" (C column) date by one day for all data rows (2-52)
$cRange = $ws.Range("C2:C52")
$cRange.Value = 46060

# Re-sort/permute rows 10-52 data (Beteckning, Datum, Markagare, Area) to match updated source export
$ws.Range("A10").Value = "A 68621-2021"
$ws.Range("B10").Value = 44502
$ws.Range("F10").ClearContents()
$ws.Range("G10").Value = 3.8

$ws.Range("A11").Value = "A 33953-2021"
$ws.Range("B11").Value = 44378
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 5

$ws.Range("A12").Value = "A 44069-2021"
$ws.Range("B12").Value = 44434
$ws.Range("F12").ClearContents()
$ws.Range("G12").Value = 1.5

$ws.Range("A13").Value = "A 9761-2021"
$ws.Range("B13").Value = 44252
$ws.Range("F13").ClearContents()
$ws.Range("G13").Value = 1.3

$ws.Range("A14").Value = "A 21264-2022"
$ws.Range("B14").Value = 44705
$ws.Range("F14").Value = "Allmännings- och besparingsskogar"
$ws.Range("G14").Value = 2.4

$ws.Range("A15").Value = "A 21972-2023"
$ws.Range("B15").Value = 45068.66232638889
$ws.Range("F15").ClearContents()
$ws.Range("G15").Value = 1.5

$ws.Range("A16").Value = "A 10263-2024"
$ws.Range("B16").Value = 45365.43090277778
$ws.Range("F16").Value = "Kyrkan"
$ws.Range("G16").Value = 1.4

$ws.Range("A17").Value = "A 30743-2021"
$ws.Range("B17").Value = 44365
$ws.Range("F17").ClearContents()
$ws.Range("G17").Value = 3

$ws.Range("A18").Value = "A 55562-2022"
$ws.Range("B18").Value = 44888
$ws.Range("F18").ClearContents()
$ws.Range("G18").Value = 0.8

$ws.Range("A19").Value = "A 27365-2025"
$ws.Range("B19").Value = 45812.64355324074
$ws.Range("F19").ClearContents()
$ws.Range("G19").Value = 11.9

$ws.Range("A20").Value = "A 50230-2024"
$ws.Range("B20").Value = 45600
$ws.Range("F20").ClearContents()
$ws.Range("G20").Value = 1.7

$ws.Range("A21").Value = "A 41546-2025"
$ws.Range("B21").Value = 45901.57927083333
$ws.Range("F21").Value = "Allmännings- och besparingsskogar"
$ws.Range("G21").Value = 4.3

$ws.Range("A22").Value = "A 41550-2025"
$ws.Range("B22").Value = 45901.58652777778
$ws.Range("F22").Value = "Allmännings- och besparingsskogar"
$ws.Range("G22").Value = 2.3

$ws.Range("A23").Value = "A 12077-2022"
$ws.Range("B23").Value = 44636.47484953704
$ws.Range("F23").ClearContents()
$ws.Range("G23").Value = 2.1

$ws.Range("A24").Value = "A 42991-2025"
$ws.Range("B24").Value = 45909.45190972222
$ws.Range("F24").ClearContents()
$ws.Range("G24").Value = 7.2

$ws.Range("A25").Value = "A 42994-2025"
$ws.Range("B25").Value = 45909.453518518516
$ws.Range("F25").ClearContents()
$ws.Range("G25").Value = 7.9

$ws.Range("A26").Value = "A 43448-2025"
$ws.Range("B26").Value = 45911.45209490741
$ws.Range("F26").ClearContents()
$ws.Range("G26").Value = 1.2

$ws.Range("A27").Value = "A 44192-2025"
$ws.Range("B27").Value = 45915.61556712963
$ws.Range("F27").ClearContents()
$ws.Range("G27").Value = 0.8

$ws.Range("A28").Value = "A 37407-2023"
$ws.Range("B28").Value = 45156.60152777778
$ws.Range("F28").ClearContents()
$ws.Range("G28").Value = 3.3

$ws.Range("A29").Value = "A 30174-2021"
$ws.Range("B29").Value = 44363
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value = 1.8

$ws.Range("A30").Value = "A 33799-2025"
$ws.Range("B30").Value = 45842.440567129626
$ws.Range("F30").Value = "Allmännings- och besparingsskogar"
$ws.Range("G30").Value = 2.5

$ws.Range("A31").Value = "A 37072-2025"
$ws.Range("B31").Value = 45875.40431712963
$ws.Range("F31").ClearContents()
$ws.Range("G31").Value = 1.5

$ws.Range("A32").Value = "A 37076-2025"
$ws.Range("B32").Value = 45875.41342592592
$ws.Range("F32").ClearContents()
$ws.Range("G32").Value = 2.4

$ws.Range("A33").Value = "A 14149-2022"
$ws.Range("B33").Value = 44651
$ws.Range("F33").Value = "Allmännings- och besparingsskogar"
$ws.Range("G33").Value = 3.8

$ws.Range("A34").Value = "A 13510-2025"
$ws.Range("B34").Value = 45736.471030092594
$ws.Range("F34").ClearContents()
$ws.Range("G34").Value = 2.5

$ws.Range("A35").Value = "A 53343-2024"
$ws.Range("B35").Value = 45614.43885416666
$ws.Range("F35").ClearContents()
$ws.Range("G35").Value = 0.9

$ws.Range("A36").Value = "A 8848-2025"
$ws.Range("B36").Value = 45713.34292824074
$ws.Range("F36").ClearContents()
$ws.Range("G36").Value = 2

$ws.Range("A37").Value = "A 58109-2025"
$ws.Range("B37").Value = 45982.597650462965
$ws.Range("F37").ClearContents()
$ws.Range("G37").Value = 2.5

$ws.Range("A38").Value = "A 58111-2025"
$ws.Range("B38").Value = 45982.59920138889
$ws.Range("F38").ClearContents()
$ws.Range("G38").Value = 0.6

$ws.Range("A39").Value = "A 58125-2025"
$ws.Range("B39").Value = 45982.615069444444
$ws.Range("F39").ClearContents()
$ws.Range("G39").Value = 0.8

$ws.Range("A40").Value = "A 54207-2025"
$ws.Range("B40").Value = 45964
$ws.Range("F40").ClearContents()
$ws.Range("G40").Value = 2.1

$ws.Range("A41").Value = "A 54203-2025"
$ws.Range("B41").Value = 45964
$ws.Range("F41").ClearContents()
$ws.Range("G41").Value = 10.3

$ws.Range("A42").Value = "A 8436-2023"
$ws.Range("B42").Value = 44977
$ws.Range("F42").Value = "Kyrkan"
$ws.Range("G42").Value = 4

$ws.Range("A43").Value = "A 1621-2026"
$ws.Range("B43").Value = 46034.47645833333
$ws.Range("F43").ClearContents()
$ws.Range("G43").Value = 1.4

$ws.Range("A44").Value = "A 1622-2026"
$ws.Range("B44").Value = 46034.47929398148
$ws.Range("F44").ClearContents()
$ws.Range("G44").Value = 1.4

$ws.Range("A45").Value = "A 7245-2025"
$ws.Range("B45").Value = 45702
$ws.Range("F45").ClearContents()
$ws.Range("G45").Value = 4

$ws.Range("A46").Value = "A 34984-2024"
$ws.Range("B46").Value = 45527
$ws.Range("F46").ClearContents()
$ws.Range("G46").Value = 4.1

$ws.Range("A47").Value = "A 62433-2025"
$ws.Range("B47").Value = 46007
$ws.Range("F47").ClearContents()
$ws.Range("G47").Value = 2

$ws.Range("A48").Value = "A 16762-2022"
$ws.Range("B48").Value = 44673.47876157407
$ws.Range("F48").ClearContents()
$ws.Range("G48").Value = 4.2

$ws.Range("A49").Value = "A 22072-2023"
$ws.Range("B49").Value = 45069
$ws.Range("F49").ClearContents()
$ws.Range("G49").Value = 3.5

$ws.Range("A50").Value = "A 30766-2022"
$ws.Range("B50").Value = 44764
$ws.Range("F50").ClearContents()
$ws.Range("G50").Value = 0.6

$ws.Range("A51").Value = "A 35036-2024"
$ws.Range("B51").Value = 45527
$ws.Range("F51").ClearContents()
$ws.Range("G51").Value = 1.7

$ws.Range("A52").Value = "A 62831-2023"
$ws.Range("B52").Value = 45270
$ws.Range("F52").ClearContents()
$ws.Range("G52").Value = 1.1
